# Update cryptos list values (prices and 1h volume % changes) per latest data pull.
# Also two rows (RenderToken / PolygonEcosystemToken) swapped ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.077.73"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.310.24"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D5").Value = "'255.20"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").Value = "'624.01"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "'1.46"
$ws.Range("E7").Value = "  +30.52%  "
$ws.Range("E8").Value = "  +5.99%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.972"
$ws.Range("E10").Value = "  +22.36%  "
$ws.Range("D11").Value = "3.306.92"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "'39.72"
$ws.Range("E13").Value = "  +11.90%  "
$ws.Range("D14").Value = "98.914.98"
$ws.Range("E14").Value = "  +1.78%  "
$ws.Range("D15").Value = "'0.0000250"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "3.932.67"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "3.306.41"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").Value = "'15.56"
$ws.Range("E20").Value = "  +3.87%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  +8.59%  "
$ws.Range("D22").Value = "'486.78"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +2.48%  "
$ws.Range("D24").Value = "'0.0000203"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'89.04"
$ws.Range("E26").Value = "  +1.13%  "
$ws.Range("D27").Value = "'11.98"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "'0.309"
$ws.Range("E28").Value = "  +29.91%  "
$ws.Range("D29").Value = "3.494.13"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D31").Value = "'0.136"
$ws.Range("E31").Value = "  +11.71%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("E33").Value = "  +11.25%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'27.87"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").Value = "'0.474"
$ws.Range("E36").Value = "  +5.38%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.22"
$ws.Range("E37").Value = "  -2.48%  "
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'24.83"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("D42").Value = "'3.64"
$ws.Range("E42").Value = "  +3.31%  "
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D46").Value = "'3.12"
$ws.Range("E46").Value = "  -5.52%  "
$ws.Range("D47").Value = "'1.95"
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").Value = "'158.37"
$ws.Range("D49").Value = "'7.33"
$ws.Range("E49").Value = "  +15.76%  "
$ws.Range("E50").Value = "  +5.68%  "
$ws.Range("D51").Value = "'4.73"
$ws.Range("E51").Value = "  +4.61%  "
